$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Add()
$ws.Name = "JsonField"
